# Update "想去人数" (F column) and "最低票价" (G column) figures that changed
# between data refreshes, on both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- 展览 sheet (rows offset by 1 vs. 全部类型 sheet for this range) ---
$ws1.Range("F4").Value = 173
$ws1.Range("F5").Value = 186
$ws1.Range("F6").Value = 2896
$ws1.Range("F8").Value = 103
$ws1.Range("F10").Value = 1603
$ws1.Range("F11").Value = 1587
$ws1.Range("F13").Value = 343
$ws1.Range("F14").Value = 234
$ws1.Range("F18").Value = 219
$ws1.Range("F23").Value = 29
$ws1.Range("F25").Value = 331
$ws1.Range("G25").Value = 52.1
$ws1.Range("F26").Value = 101
$ws1.Range("G26").Value = 55
$ws1.Range("F27").Value = 88
$ws1.Range("F28").Value = 8
$ws1.Range("F29").Value = 1892
$ws1.Range("F31").Value = 439
$ws1.Range("F33").Value = 131
$ws1.Range("F34").Value = 573
$ws1.Range("F36").Value = 322
$ws1.Range("F38").Value = 473

# --- 全部类型 sheet ---
$ws4.Range("F5").Value = 173
$ws4.Range("F6").Value = 186
$ws4.Range("F7").Value = 2896
$ws4.Range("F9").Value = 103
$ws4.Range("F11").Value = 1603
$ws4.Range("F12").Value = 1587
$ws4.Range("F14").Value = 343
$ws4.Range("F15").Value = 234
$ws4.Range("F19").Value = 219
$ws4.Range("F24").Value = 29
$ws4.Range("F26").Value = 331
$ws4.Range("G26").Value = 52.1
$ws4.Range("F27").Value = 101
$ws4.Range("G27").Value = 55
$ws4.Range("F28").Value = 88
$ws4.Range("F29").Value = 8
$ws4.Range("F30").Value = 1892
$ws4.Range("F32").Value = 439
$ws4.Range("F34").Value = 131
$ws4.Range("F35").Value = 573
$ws4.Range("F37").Value = 322
$ws4.Range("F39").Value = 473

$wb.Save()
